# bergupdate: almost done with ASCAT slides.
#
# The deck's "Update automatically" date footer (the datetimeFigureOut
# field cached on the slide master and on every slide layout) was
# refreshed from 4/14/2022 to 4/15/2022. Walk the slide master plus all
# of its custom layouts and, for every shape whose current text is the
# old cached date, replace it with the new one.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster

$oldDate = "4/14/2022"
$newDate = "4/15/2022"

function Update-DateShape($sh) {
    if (-not $sh.HasTextFrame) {
        return
    }
    $tf = $sh.TextFrame
    if (-not $tf.HasText) {
        return
    }
    $tr = $tf.TextRange
    if ($tr.Text -eq $oldDate) {
        $tr.Text = $newDate
    }
}

# Slide master's own Date placeholder.
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateShape $master.Shapes.Item($i)
}

# Every slide layout's Date placeholder (each layout has its own copy of
# the footer placeholders, and the cached date text lives on each one).
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        Update-DateShape $layout.Shapes.Item($i)
    }
}
